# Apply the edit described by the diff:
# - Insert a new row 2 (Chilean Primera B, Deportes Concepcion vs Antofagasta)
#   which pushes the existing Brazilian Serie A row down to row 3.
# - Tweak a few odds values on the (now) row 3.
# - Append two brand-new match rows (4 and 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row above the current row 2 (shifts Brazilian Serie A -> row 3)
#    Use xlFormatFromRightOrBelow (1) so the new row does not inherit the bold
#    header formatting from row 1 above it.
$ws.Rows.Item(2).Insert([Type]::Missing, 1)
$ws.Rows.Item(2).ClearFormats()

# 2. Populate the new row 2 with the Chilean Primera B (18:00) match
#    Dates/times are stored as plain text in the source file, so force the
#    cell format to Text before writing them to avoid Excel auto-converting
#    them into date/time serial numbers.
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 3).NumberFormat = "@"

$row2 = @("Chilean Primera B","2025-11-18","18:00:00","Deportes Concepcion","Antofagasta", `
    1.93,2.62,3.1,5,2.88,6.4,1.01,1.01,2.32,1.01,1.62,1.94,1.19,2.6,1.61,1.69,1.25,1.61, `
    1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

$ws.Cells.Item(2, 2).NumberFormat = "General"
$ws.Cells.Item(2, 3).NumberFormat = "General"

# 3. Fix up the odds that changed on the Brazilian Serie A row, now row 3
$ws.Cells.Item(3, 7).Value = 1.31    # G3 Odd_H_Lay
$ws.Cells.Item(3, 15).Value = 1.22   # O3 Odd_Over15_FT_Back
$ws.Cells.Item(3, 23).Value = 4.2    # W3 Double_Chance_Draw_or_Away_Back

# 4. Add the Chilean Primera B (20:30) match as row 4
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 3).NumberFormat = "@"

$row4 = @("Chilean Primera B","2025-11-18","20:30:00","CSD Rangers","San Marcos", `
    1.65,2.1,1.91,1000,3.15,1000,0,0,0,0,1.6,1.94,0,0,0,0,0,0, `
    0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}

# 5. Add the Colombian Primera A match as row 5
$ws.Cells.Item(5, 2).NumberFormat = "@"
$ws.Cells.Item(5, 3).NumberFormat = "@"

$row5 = @("Colombian Primera A","2025-11-18","20:30:00","Fortaleza FC","Tolima", `
    1.98,2.2,4,5.5,3.15,3.7,0,0,0,0,1.66,2.2,0,0,0,0,0,0, `
    0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5[$i]
}
